# QB siteyml change 1/21
# SocialSupport.xlsx: replace the "Occurrence" column (E) with a "Dates Used"
# column that lists the actual survey date ranges instead of question numbers,
# and update the sheet selection to the full column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header: E1 "Occurrence" -> "Dates Used"
$ws.Range("E1").Value = "Dates Used"

# New "Dates Used" values per row, replacing the old occurrence numbers
$ws.Range("E2").Value = "July 13 - July 16, Aug 10 - Aug 13, Aug 24 - Aug 27, Sept 8 - Sept 10, Dec 1 - Dec 3, Feb 3 - Feb 6, April 14 - April 16, April 28 - May 1, May 25 - May 28, July 21 - July 23, September 1 - September 3, January 11 - January 18"
$ws.Range("E3").Value = "July 13 - July 16, Aug 10 - Aug 13, Aug 24 - Aug 27, Sept 8 - Sept 10, Dec 1 - Dec 3, Feb 3 - Feb 6, April 14 - April 16, April 28 - May 1, May 25 - May 28, July 21 - July 23, September 1 - September 3, October 27 - November 2"
$ws.Range("E4").Value = "July 13 - July 16"
$ws.Range("E5").Value = "July 13 - July 16"
$ws.Range("E6").Value = "Aug 10 - Aug 13, Aug 24 - Aug 27"
$ws.Range("E7").Value = "Aug 10 - Aug 13, Aug 24 - Aug 27"
$ws.Range("E8").Value = "Aug 10 - Aug 13, Aug 24 - Aug 27"

# Match the updated view selection (whole column E selected)
$ws.Range("E1:E1048576").Select()
